# Sumit's resume: add the GitHub profile URL to the placeholder line that
# previously held only tab characters, and bump that line's font size.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item(6)          # "Text Placeholder 26"

$tf = $shp.TextFrame
$tr = $tf.TextRange

# The 7th paragraph in this placeholder is a standalone run of tab
# characters ("\t\t\t") that sits between "Check out my work on GitHub"
# and a trailing blank line. Grab that single run and update it in place.
$para = $tr.Paragraphs(7, 1)
$run = $para.Runs(1, 1)

$run.Text = "`t`thttps://github.com/SumitGangotri`t"
$run.Font.Size = 14
